# Switch the deck's theme colors from the custom "Integral" / "Red Violet"
# palette over to the built-in "Office Theme" palette (the same colors
# that ship as the Office color-scheme default).
#
# PowerPoint's Theme Colors gallery writes all twelve slots (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) of the active ThemeColorScheme;
# reproduce that here one slot at a time.

function RgbHexToCom([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = RgbHexToCom $officeColors[$i - 1]
}
